$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": update PORCELANATO sale for BECERRA FARIAS ROSA DAYANA (row 3)
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M3").Value = 6263.28

# Sheet "VENTA MENSUAL": update julio sale for BECERRA FARIAS ROSA DAYANA (row 3) and its total (row 24)
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F3").Value = 6836.54
$wsVentaMensual.Range("F24").Value = 57458

# Sheet "CUMPLIMIENTO MENSUAL": update PORCELANATO group totals (row 16) and grand total (row 19)
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 52125.18
$wsCumplimiento.Range("E16").Value = -13368.64
$wsCumplimiento.Range("F16").Value = 1.344938944498141

$wsCumplimiento.Range("D19").Value = 57458
$wsCumplimiento.Range("E19").Value = 765.0038630460349
$wsCumplimiento.Range("F19").Value = 0.9868607970683633
